# Update the "想去人数" (number of people wanting to attend) figures that
# changed when the page data was regenerated (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F7").Value  = 2703
$wsExhibit.Range("F11").Value = 10197
$wsExhibit.Range("F16").Value = 11808
$wsExhibit.Range("F17").Value = 12214

# Sheet "全部类型" (All types) - same events, shifted by one row because of
# an extra entry earlier in the sheet.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value  = 2703
$wsAll.Range("F12").Value = 10197
$wsAll.Range("F17").Value = 11808
$wsAll.Range("F18").Value = 12214
